$wb = $excel.ActiveWorkbook
$wsInput = $wb.Worksheets.Item("Input")
$wsOutput = $wb.Worksheets.Item("Output")

# Update the client name used for AddClientMember step
$wsInput.Range("B7").Value = "Jhon Deer"

# Replace the old "verify1 / Client Not Attached" step with the new
# "GroupAddClient / click" step
$wsInput.Range("A8").Value = "GroupAddClient"
$wsInput.Range("B8").Value = "click"

# Add the new verification row on the Output sheet
$wsOutput.Range("A2").Value = "VerifyClientCreated"
$wsOutput.Range("B2").Value = "Jhon Deer"

# Update selections / active sheet to match the authored state
$wsInput.Range("A12").Select()
$wsOutput.Activate()
$wsOutput.Range("A6").Select()
